# Scheduled runner update: refresh market-derived profit metrics
# (currentAveragePrice/NQ/HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# for a handful of leves across the ALC, ARM, BSM, CRP, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 5007293.5
$ws.Range("I80").Value = 24280.4
$ws.Range("J80").Value = 6253047
$ws.Range("K80").Value = 72841.20000000001
$ws.Range("L80").Value = 18759141
$ws.Range("M80").Value = -71843.20000000001
$ws.Range("N80").Value = -18761137

$ws.Range("H83").Value = 5007293.5
$ws.Range("I83").Value = 24280.4
$ws.Range("J83").Value = 6253047
$ws.Range("K83").Value = 218523.6
$ws.Range("L83").Value = 56277423
$ws.Range("M83").Value = -213531.6
$ws.Range("N83").Value = -56287407

$ws.Range("H98").Value = 20329392
$ws.Range("I98").Value = 6452596.5
$ws.Range("J98").Value = 127874570
$ws.Range("K98").Value = 6452596.5
$ws.Range("L98").Value = 127874570
$ws.Range("M98").Value = -6451098.5
$ws.Range("N98").Value = -127877566

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H107").Value = 871.25
$ws.Range("I107").Value = 974.5238000000001
$ws.Range("J107").Value = 148.33333
$ws.Range("K107").Value = 974.5238000000001
$ws.Range("L107").Value = 148.33333
$ws.Range("M107").Value = 945.4761999999999
$ws.Range("N107").Value = -3988.33333

$ws.Range("H113").Value = 5884045.5
$ws.Range("I113").Value = 8334848.5
$ws.Range("J113").Value = 2119.2
$ws.Range("K113").Value = 8334848.5
$ws.Range("L113").Value = 2119.2
$ws.Range("M113").Value = -8331594.5
$ws.Range("N113").Value = -8627.200000000001

$ws.Range("H116").Value = 11722468
$ws.Range("I116").Value = 5130513.5
$ws.Range("J116").Value = 23964670
$ws.Range("K116").Value = 5130513.5
$ws.Range("L116").Value = 23964670
$ws.Range("M116").Value = -5127071.5
$ws.Range("N116").Value = -23971554

$ws.Range("H122").Value = 20329392
$ws.Range("I122").Value = 6452596.5
$ws.Range("J122").Value = 127874570
$ws.Range("K122").Value = 19357789.5
$ws.Range("L122").Value = 383623710
$ws.Range("M122").Value = -19355339.5
$ws.Range("N122").Value = -383628610

$ws.Range("H132").Value = 1425781.5
$ws.Range("I132").Value = 1019.6377
$ws.Range("J132").Value = 12348956
$ws.Range("K132").Value = 3058.9131
$ws.Range("L132").Value = 37046868
$ws.Range("M132").Value = -528.9130999999998
$ws.Range("N132").Value = -37051928

$ws.Range("H135").Value = 12820912
$ws.Range("I135").Value = 376.69446
$ws.Range("J135").Value = 166667330
$ws.Range("K135").Value = 3390.25014
$ws.Range("L135").Value = 1500005970
$ws.Range("M135").Value = -855.2501400000001
$ws.Range("N135").Value = -1500011040

$ws.Range("H137").Value = 10974424
$ws.Range("I137").Value = 831.46344
$ws.Range("K137").Value = 2494.39032
$ws.Range("M137").Value = 55.60968000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6584696
$ws.Range("I32").Value = 1459134.6
$ws.Range("J32").Value = 55562280
$ws.Range("K32").Value = 1459134.6
$ws.Range("L32").Value = 55562280
$ws.Range("M32").Value = -1458847.6
$ws.Range("N32").Value = -55562854

$ws.Range("H61").Value = 1570902.9
$ws.Range("I61").Value = 786795.9
$ws.Range("J61").Value = 5348873.5
$ws.Range("K61").Value = 786795.9
$ws.Range("L61").Value = 5348873.5
$ws.Range("M61").Value = -786583.9
$ws.Range("N61").Value = -5349297.5

$ws.Range("H136").Value = 1570902.9
$ws.Range("I136").Value = 786795.9
$ws.Range("J136").Value = 5348873.5
$ws.Range("K136").Value = 2360387.7
$ws.Range("L136").Value = 16046620.5
$ws.Range("M136").Value = -2357837.7
$ws.Range("N136").Value = -16051720.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1563.2632
$ws.Range("I99").Value = 1336.6364
$ws.Range("J99").Value = 1874.875
$ws.Range("K99").Value = 1336.6364
$ws.Range("L99").Value = 1874.875
$ws.Range("M99").Value = 161.3635999999999
$ws.Range("N99").Value = -4870.875

$ws.Range("H134").Value = 8371285.5
$ws.Range("I134").Value = 9091621
$ws.Range("J134").Value = 3969231.5
$ws.Range("K134").Value = 27274863
$ws.Range("L134").Value = 11907694.5
$ws.Range("M134").Value = -27272328
$ws.Range("N134").Value = -11912764.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1080.909
$ws.Range("I16").Value = 1122.2222
$ws.Range("J16").Value = 895
$ws.Range("K16").Value = 1122.2222
$ws.Range("L16").Value = 895
$ws.Range("M16").Value = -835.2221999999999
$ws.Range("N16").Value = -1469

$ws.Range("H31").Value = 1527970.4
$ws.Range("I31").Value = 1662.0385
$ws.Range("J31").Value = 4173571.5
$ws.Range("K31").Value = 1662.0385
$ws.Range("L31").Value = 4173571.5
$ws.Range("M31").Value = -1367.0385
$ws.Range("N31").Value = -4174161.5

$ws.Range("H34").Value = 1527970.4
$ws.Range("I34").Value = 1662.0385
$ws.Range("J34").Value = 4173571.5
$ws.Range("K34").Value = 1662.0385
$ws.Range("L34").Value = 4173571.5
$ws.Range("M34").Value = -1460.0385
$ws.Range("N34").Value = -4173975.5

$ws.Range("H113").Value = 1080.909
$ws.Range("I113").Value = 1122.2222
$ws.Range("J113").Value = 895
$ws.Range("K113").Value = 1122.2222
$ws.Range("L113").Value = 895
$ws.Range("M113").Value = 1047.7778
$ws.Range("N113").Value = -5235

$ws.Range("H132").Value = 2402.889
$ws.Range("I132").Value = 1374.2142
$ws.Range("K132").Value = 4122.642599999999
$ws.Range("M132").Value = -1592.642599999999

$ws.Range("H134").Value = 854750.6
$ws.Range("I134").Value = 3527.3333
$ws.Range("J134").Value = 5004464
$ws.Range("K134").Value = 10581.9999
$ws.Range("L134").Value = 15013392
$ws.Range("M134").Value = -8046.999899999999
$ws.Range("N134").Value = -15018462

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3585933.8
$ws.Range("I136").Value = 5849443
$ws.Range("J136").Value = 2044.1666
$ws.Range("K136").Value = 17548329
$ws.Range("L136").Value = 6132.4998
$ws.Range("M136").Value = -17545779
$ws.Range("N136").Value = -11232.4998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5232.6
$ws.Range("I136").Value = 3120.4595
$ws.Range("J136").Value = 15001.25
$ws.Range("K136").Value = 9361.378499999999
$ws.Range("L136").Value = 45003.75
$ws.Range("M136").Value = -6811.378499999999
$ws.Range("N136").Value = -50103.75

